$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits inside the "Ethernet phy (...)"
#    paragraph (between "figure 92" and " STM43F407i datasheet)"). In the
#    target document that bookmark has moved into the new ".1in header "
#    bullet further down. Delete the old one now; it will be re-created at
#    the right spot once that paragraph exists.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Locate the anchor paragraphs by their (unique, stable) text.
# ---------------------------------------------------------------------------
function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -eq ($text + [char]13)) {
            return $i
        }
    }
    return -1
}

$edgeIdx = Find-ParaIndex "Edge connectors"

# ---------------------------------------------------------------------------
# 3. Insert the new "Edge connectors" sub-bullets (".1in header " bookmark,
#    " 2mm (waveshare)") plus the new "Programming connector" group right
#    after "Edge connectors" and before "USB FS connector".
# ---------------------------------------------------------------------------
$idx = $edgeIdx

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = ".1in header "
$d.Paragraphs($idx).Range.ListFormat.ListLevelNumber = 2

# Re-create "_GoBack" right after the text we just typed (collapsed, same
# shape as in the source document).
$p = $d.Paragraphs($idx)
$bmStart = $p.Range.Start + ".1in header ".Length
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = " 2mm (waveshare)"
$d.Paragraphs($idx).Range.ListFormat.ListLevelNumber = 2

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "Programming connector "
$d.Paragraphs($idx).Range.ListFormat.ListLevelNumber = 1

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = ".05 in 10pin SWD connector "
$d.Paragraphs($idx).Range.ListFormat.ListLevelNumber = 2

$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.Text = "keyed shroud"
$d.Paragraphs($idx).Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------------
# 4. "USB FS connector" keeps its text; replace its two sub-bullets
#    ("Mini?" / "Macro?") with the single "Micro" sub-bullet.
# ---------------------------------------------------------------------------
$usbIdx = Find-ParaIndex "USB FS connector"
$miniIdx = Find-ParaIndex "Mini?"
$macroIdx = Find-ParaIndex "Macro?"

# Delete "Macro?" first so "Mini?"'s index stays valid.
$d.Paragraphs($macroIdx).Range.Delete()
$d.Paragraphs($miniIdx).Range.Text = "Micro"

# ---------------------------------------------------------------------------
# 5. Add the new sub-bullet under "Power input barrel jack ".
# ---------------------------------------------------------------------------
$powerIdx = Find-ParaIndex "Power input barrel jack "
$d.Paragraphs($powerIdx).Range.InsertParagraphAfter()
$newIdx = $powerIdx + 1
$d.Paragraphs($newIdx).Range.Text = "5.1 x 2.0 mm ish"
$d.Paragraphs($newIdx).Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------------
# 6. Add the new sub-bullet under "SDIO memory card interface?".
# ---------------------------------------------------------------------------
$sdioIdx = Find-ParaIndex "SDIO memory card interface?"
$d.Paragraphs($sdioIdx).Range.InsertParagraphAfter()
$newIdx = $sdioIdx + 1
$d.Paragraphs($newIdx).Range.Text = "Friction lock (push pull)"
$d.Paragraphs($newIdx).Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------------
# 7. Append the new "System status LEDs" block at the end of the list.
# ---------------------------------------------------------------------------
$lastIdx = Find-ParaIndex "Friction lock (push pull)"

$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "System status LEDs "
$d.Paragraphs($lastIdx).Range.ListFormat.ListLevelNumber = 1

$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "Count 8 LED" + [char]8217 + "s (code for marionette state) - mBus"
$d.Paragraphs($lastIdx).Range.ListFormat.ListLevelNumber = 2

$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "RGB (heartbeat)"
$d.Paragraphs($lastIdx).Range.ListFormat.ListLevelNumber = 2

Write-Output "done"
